$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-30 12:55:45"

# Update the timestamp column (O) for all data rows (2 through 397)
for ($r = 2; $r -le 397; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Row 222: ratingAmount (D) changes from 11 to 12
$ws.Cells.Item(222, 4).Value = 12

# Row 313: ratingAmount (D) changes from 3 to 4, ratingValue (E) changes from 4 to 3.5
$ws.Cells.Item(313, 4).Value = 4
$ws.Cells.Item(313, 5).Value = 3.5
